# The generated site dropped the trailing "Ver no Jupiter / Salvar em pdf /
# Salvar em docx" line, the footer copyright notice, and the blank paragraph
# that separated them from the requirements list. Remove that whole block,
# keeping the last requirement line ("LOQ4073: ...") intact and leaving the
# blank paragraph (and page-break paragraph) that originally followed the
# footer untouched.
$d = $word.ActiveDocument

# Locate the end of the last requirement paragraph that must be preserved.
$anchor = $d.Content
$anchor.Find.Execute("LOQ4073: Química Geral II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor.MoveEnd(1, 1) | Out-Null   # include its paragraph mark

# Locate the end of the copyright paragraph that must be removed.
$footerEnd = $d.Content
$footerEnd.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$footerEnd.MoveEnd(1, 1) | Out-Null   # include its paragraph mark

# Delete everything from right after the requirement paragraph's mark
# through the end of the copyright paragraph's mark: the blank paragraph,
# the "Ver no Jupiter..." paragraph, and the copyright paragraph.
$d.Range($anchor.End, $footerEnd.End).Delete()
